$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 421.69232  # H8: 481.75 -> 421.69232
$ws.Cells.Item(8, 9).Value = 563.6667  # I8: 845.25 -> 563.6667
$ws.Cells.Item(8, 11).Value = 1691.0001  # K8: 2535.75 -> 1691.0001
$ws.Cells.Item(8, 13).Value = -1552.0001  # M8: -2396.75 -> -1552.0001

$ws.Cells.Item(100, 8).Value = 2972.182  # H100: 2993.625 -> 2972.182
$ws.Cells.Item(100, 9).Value = 2969.4  # I100: 2992.7144 -> 2969.4
$ws.Cells.Item(100, 11).Value = 2969.4  # K100: 2992.7144 -> 2969.4
$ws.Cells.Item(100, 13).Value = -2428.4  # M100: -2451.7144 -> -2428.4

$ws.Cells.Item(113, 8).Value = 4005  # H113: 2000 -> 4005
$ws.Cells.Item(113, 9).Value = 4005  # I113: 0 -> 4005
$ws.Cells.Item(113, 10).Value = 0  # J113: 2000 -> 0
$ws.Cells.Item(113, 11).Value = 4005  # K113: 0 -> 4005
$ws.Cells.Item(113, 12).Value = 0  # L113: 2000 -> 0
$ws.Cells.Item(113, 13).Value = -751  # M113: None -> -751
$ws.Cells.Item(113, 14).ClearContents()  # N113: -8508 -> (removed)

$ws.Cells.Item(116, 8).Value = 3312.7778  # H116: 3489.375 -> 3312.7778
$ws.Cells.Item(116, 9).Value = 1940  # I116: 1980 -> 1940
$ws.Cells.Item(116, 11).Value = 1940  # K116: 1980 -> 1940
$ws.Cells.Item(116, 13).Value = 1502  # M116: 1462 -> 1502

$ws.Cells.Item(135, 8).Value = 3488  # H135: 3335.3635 -> 3488
$ws.Cells.Item(135, 9).Value = 640.6667  # I135: 807.5714 -> 640.6667
$ws.Cells.Item(135, 11).Value = 5766.0003  # K135: 7268.1426 -> 5766.0003
$ws.Cells.Item(135, 13).Value = -3231.0003  # M135: -4733.1426 -> -3231.0003

$ws.Cells.Item(137, 8).Value = 1360.9565  # H137: 1465.8718 -> 1360.9565
$ws.Cells.Item(137, 9).Value = 935.5  # I137: 976.86365 -> 935.5
$ws.Cells.Item(137, 10).Value = 1914.05  # J137: 2098.7058 -> 1914.05
$ws.Cells.Item(137, 11).Value = 2806.5  # K137: 2930.59095 -> 2806.5
$ws.Cells.Item(137, 12).Value = 5742.15  # L137: 6296.117400000001 -> 5742.15
$ws.Cells.Item(137, 13).Value = -256.5  # M137: -380.5909499999998 -> -256.5
$ws.Cells.Item(137, 14).Value = -10842.15  # N137: -11396.1174 -> -10842.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1400  # H2: 9092.416999999999 -> 1400
$ws.Cells.Item(2, 9).Value = 1200  # I2: 11735.444 -> 1200
$ws.Cells.Item(2, 10).Value = 1600  # J2: 1163.3334 -> 1600
$ws.Cells.Item(2, 11).Value = 1200  # K2: 11735.444 -> 1200
$ws.Cells.Item(2, 12).Value = 1600  # L2: 1163.3334 -> 1600
$ws.Cells.Item(2, 13).Value = -1087  # M2: -11622.444 -> -1087
$ws.Cells.Item(2, 14).Value = -1826  # N2: -1389.3334 -> -1826

$ws.Cells.Item(63, 8).Value = 2134.2354  # H63: 2136.7354 -> 2134.2354
$ws.Cells.Item(63, 9).Value = 2043.6522  # I63: 2047.3478 -> 2043.6522
$ws.Cells.Item(63, 11).Value = 2043.6522  # K63: 2047.3478 -> 2043.6522
$ws.Cells.Item(63, 13).Value = -1357.6522  # M63: -1361.3478 -> -1357.6522

$ws.Cells.Item(66, 8).Value = 2134.2354  # H66: 2136.7354 -> 2134.2354
$ws.Cells.Item(66, 9).Value = 2043.6522  # I66: 2047.3478 -> 2043.6522
$ws.Cells.Item(66, 11).Value = 10218.261  # K66: 10236.739 -> 10218.261
$ws.Cells.Item(66, 13).Value = -6786.261  # M66: -6804.739 -> -6786.261

$ws.Cells.Item(74, 8).Value = 1461.7894  # H74: 1647.2 -> 1461.7894
$ws.Cells.Item(74, 9).Value = 1073.375  # I74: 1316 -> 1073.375
$ws.Cells.Item(74, 10).Value = 3533.3333  # J74: 3800 -> 3533.3333
$ws.Cells.Item(74, 11).Value = 1073.375  # K74: 1316 -> 1073.375
$ws.Cells.Item(74, 12).Value = 3533.3333  # L74: 3800 -> 3533.3333
$ws.Cells.Item(74, 13).Value = -199.375  # M74: -442 -> -199.375
$ws.Cells.Item(74, 14).Value = -5281.3333  # N74: -5548 -> -5281.3333

$ws.Cells.Item(77, 8).Value = 1461.7894  # H77: 1647.2 -> 1461.7894
$ws.Cells.Item(77, 9).Value = 1073.375  # I77: 1316 -> 1073.375
$ws.Cells.Item(77, 10).Value = 3533.3333  # J77: 3800 -> 3533.3333
$ws.Cells.Item(77, 11).Value = 5366.875  # K77: 6580 -> 5366.875
$ws.Cells.Item(77, 12).Value = 17666.6665  # L77: 19000 -> 17666.6665
$ws.Cells.Item(77, 13).Value = -998.875  # M77: -2212 -> -998.875
$ws.Cells.Item(77, 14).Value = -26402.6665  # N77: -27736 -> -26402.6665

$ws.Cells.Item(102, 8).Value = 33335290  # H102: 41668660 -> 33335290
$ws.Cells.Item(102, 9).Value = 55556830  # I102: 83334340 -> 55556830
$ws.Cells.Item(102, 11).Value = 55556830  # K102: 83334340 -> 55556830
$ws.Cells.Item(102, 13).Value = -55555208  # M102: -83332718 -> -55555208

$ws.Cells.Item(107, 8).Value = 0  # H107: 30000 -> 0
$ws.Cells.Item(107, 10).Value = 0  # J107: 30000 -> 0
$ws.Cells.Item(107, 12).Value = 0  # L107: 30000 -> 0
$ws.Cells.Item(107, 14).ClearContents()  # N107: -37680 -> (removed)

$ws.Cells.Item(116, 8).Value = 1400  # H116: 9092.416999999999 -> 1400
$ws.Cells.Item(116, 9).Value = 1200  # I116: 11735.444 -> 1200
$ws.Cells.Item(116, 10).Value = 1600  # J116: 1163.3334 -> 1600
$ws.Cells.Item(116, 11).Value = 1200  # K116: 11735.444 -> 1200
$ws.Cells.Item(116, 12).Value = 1600  # L116: 1163.3334 -> 1600
$ws.Cells.Item(116, 13).Value = 1094  # M116: -9441.444 -> 1094
$ws.Cells.Item(116, 14).Value = -6188  # N116: -5751.3334 -> -6188

$ws.Cells.Item(122, 8).Value = 857.5294  # H122: 998.73334 -> 857.5294
$ws.Cells.Item(122, 9).Value = 938.5333000000001  # I122: 1113.9231 -> 938.5333000000001
$ws.Cells.Item(122, 11).Value = 2815.5999  # K122: 3341.7693 -> 2815.5999
$ws.Cells.Item(122, 13).Value = -365.5999000000002  # M122: -891.7692999999999 -> -365.5999000000002

$ws.Cells.Item(132, 8).Value = 2823.9473  # H132: 1741.1389 -> 2823.9473
$ws.Cells.Item(132, 9).Value = 2410.4666  # I132: 1393 -> 2410.4666
$ws.Cells.Item(132, 10).Value = 4374.5  # J132: 3899.6 -> 4374.5
$ws.Cells.Item(132, 11).Value = 7231.399800000001  # K132: 4179 -> 7231.399800000001
$ws.Cells.Item(132, 12).Value = 13123.5  # L132: 11698.8 -> 13123.5
$ws.Cells.Item(132, 13).Value = -4701.399800000001  # M132: -1649 -> -4701.399800000001
$ws.Cells.Item(132, 14).Value = -18183.5  # N132: -16758.8 -> -18183.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1400  # H3: 9092.416999999999 -> 1400
$ws.Cells.Item(3, 9).Value = 1200  # I3: 11735.444 -> 1200
$ws.Cells.Item(3, 10).Value = 1600  # J3: 1163.3334 -> 1600
$ws.Cells.Item(3, 11).Value = 1200  # K3: 11735.444 -> 1200
$ws.Cells.Item(3, 12).Value = 1600  # L3: 1163.3334 -> 1600
$ws.Cells.Item(3, 13).Value = -1086  # M3: -11621.444 -> -1086
$ws.Cells.Item(3, 14).Value = -1828  # N3: -1391.3334 -> -1828

$ws.Cells.Item(134, 8).Value = 7329.9473  # H134: 8037.9414 -> 7329.9473
$ws.Cells.Item(134, 9).Value = 1559.4615  # I134: 1604.4546 -> 1559.4615
$ws.Cells.Item(134, 11).Value = 4678.3845  # K134: 4813.3638 -> 4678.3845
$ws.Cells.Item(134, 13).Value = -2143.3845  # M134: -2278.3638 -> -2143.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 100001330  # H16: 55556772 -> 100001330
$ws.Cells.Item(16, 9).Value = 142858450  # I16: 90910370 -> 142858450
$ws.Cells.Item(16, 10).Value = 1366.6666  # J16: 1122.8572 -> 1366.6666
$ws.Cells.Item(16, 11).Value = 142858450  # K16: 90910370 -> 142858450
$ws.Cells.Item(16, 12).Value = 1366.6666  # L16: 1122.8572 -> 1366.6666
$ws.Cells.Item(16, 13).Value = -142858163  # M16: -90910083 -> -142858163
$ws.Cells.Item(16, 14).Value = -1940.6666  # N16: -1696.8572 -> -1940.6666

$ws.Cells.Item(113, 8).Value = 100001330  # H113: 55556772 -> 100001330
$ws.Cells.Item(113, 9).Value = 142858450  # I113: 90910370 -> 142858450
$ws.Cells.Item(113, 10).Value = 1366.6666  # J113: 1122.8572 -> 1366.6666
$ws.Cells.Item(113, 11).Value = 142858450  # K113: 90910370 -> 142858450
$ws.Cells.Item(113, 12).Value = 1366.6666  # L113: 1122.8572 -> 1366.6666
$ws.Cells.Item(113, 13).Value = -142856280  # M113: -90908200 -> -142856280
$ws.Cells.Item(113, 14).Value = -5706.6666  # N113: -5462.8572 -> -5706.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 3176.5334  # H39: 3303.4285 -> 3176.5334
$ws.Cells.Item(39, 10).Value = 3132  # J39: 3265.2307 -> 3132
$ws.Cells.Item(39, 12).Value = 9396  # L39: 9795.6921 -> 9396
$ws.Cells.Item(39, 14).Value = -9984  # N39: -10383.6921 -> -9984

$ws.Cells.Item(51, 8).Value = 1357.1428  # H51: 201890.8 -> 1357.1428
$ws.Cells.Item(51, 9).Value = 1333.3334  # I51: 334484.66 -> 1333.3334
$ws.Cells.Item(51, 10).Value = 1500  # J51: 3000 -> 1500
$ws.Cells.Item(51, 11).Value = 4000.0002  # K51: 1003453.98 -> 4000.0002
$ws.Cells.Item(51, 12).Value = 4500  # L51: 9000 -> 4500
$ws.Cells.Item(51, 13).Value = -3540.0002  # M51: -1002993.98 -> -3540.0002
$ws.Cells.Item(51, 14).Value = -5420  # N51: -9920 -> -5420

$ws.Cells.Item(55, 8).Value = 2063.6365  # H55: 2220 -> 2063.6365

$ws.Cells.Item(109, 8).Value = 143947.58  # H109: 144968.72 -> 143947.58
$ws.Cells.Item(109, 9).Value = 167372.17  # I109: 251145.25 -> 167372.17
$ws.Cells.Item(109, 11).Value = 502116.51  # K109: 753435.75 -> 502116.51
$ws.Cells.Item(109, 13).Value = -501076.51  # M109: -752395.75 -> -501076.51

$ws.Cells.Item(113, 8).Value = 724  # H113: 653.75 -> 724
$ws.Cells.Item(113, 9).Value = 0  # I113: 499.2 -> 0
$ws.Cells.Item(113, 11).Value = 0  # K113: 1497.6 -> 0
$ws.Cells.Item(113, 13).ClearContents()  # M113: 672.4000000000001 -> (removed)

$ws.Cells.Item(118, 8).Value = 865.8  # H118: 769.5 -> 865.8
$ws.Cells.Item(118, 9).Value = 582.25  # I118: 523.4 -> 582.25
$ws.Cells.Item(118, 11).Value = 1746.75  # K118: 1570.2 -> 1746.75
$ws.Cells.Item(118, 13).Value = -503.75  # M118: -327.1999999999998 -> -503.75

$ws.Cells.Item(121, 8).Value = 590.1111  # H121: 752.55554 -> 590.1111
$ws.Cells.Item(121, 9).Value = 406.33334  # I121: 450 -> 406.33334
$ws.Cells.Item(121, 10).Value = 957.6667  # J121: 994.6 -> 957.6667
$ws.Cells.Item(121, 11).Value = 1219.00002  # K121: 1350 -> 1219.00002
$ws.Cells.Item(121, 12).Value = 2873.0001  # L121: 2983.8 -> 2873.0001
$ws.Cells.Item(121, 13).Value = 90.99998000000005  # M121: -40 -> 90.99998000000005
$ws.Cells.Item(121, 14).Value = -5493.0001  # N121: -5603.8 -> -5493.0001

$ws.Cells.Item(122, 8).Value = 960  # H122: 1106.421 -> 960
$ws.Cells.Item(122, 9).Value = 526.5714  # I122: 849.3333 -> 526.5714
$ws.Cells.Item(122, 10).Value = 1128.5555  # J122: 1154.625 -> 1128.5555
$ws.Cells.Item(122, 11).Value = 4739.1426  # K122: 7643.9997 -> 4739.1426
$ws.Cells.Item(122, 12).Value = 10156.9995  # L122: 10391.625 -> 10156.9995
$ws.Cells.Item(122, 13).Value = -2289.1426  # M122: -5193.9997 -> -2289.1426
$ws.Cells.Item(122, 14).Value = -15056.9995  # N122: -15291.625 -> -15056.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2234.5454  # H126: 2132.7273 -> 2234.5454
$ws.Cells.Item(126, 9).Value = 1825.7142  # I126: 1820 -> 1825.7142
$ws.Cells.Item(126, 10).Value = 2950  # J126: 2966.6667 -> 2950
$ws.Cells.Item(126, 11).Value = 5477.142599999999  # K126: 5460 -> 5477.142599999999
$ws.Cells.Item(126, 12).Value = 8850  # L126: 8900.000100000001 -> 8850
$ws.Cells.Item(126, 13).Value = -3007.142599999999  # M126: -2990 -> -3007.142599999999
$ws.Cells.Item(126, 14).Value = -13790  # N126: -13840.0001 -> -13790

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1137.125  # H16: 746.3077 -> 1137.125
$ws.Cells.Item(16, 9).Value = 1192.9333  # I16: 746.3077 -> 1192.9333
$ws.Cells.Item(16, 10).Value = 300  # J16: 0 -> 300
$ws.Cells.Item(16, 11).Value = 1192.9333  # K16: 746.3077 -> 1192.9333
$ws.Cells.Item(16, 12).Value = 300  # L16: 0 -> 300
$ws.Cells.Item(16, 13).Value = -1022.9333  # M16: -576.3077 -> -1022.9333
$ws.Cells.Item(16, 14).Value = -640  # N16: None -> -640

$ws.Cells.Item(26, 8).Value = 2504.5  # H26: 2502.25 -> 2504.5
$ws.Cells.Item(26, 9).Value = 2504.5  # I26: 2503 -> 2504.5
$ws.Cells.Item(26, 10).Value = 0  # J26: 2500 -> 0
$ws.Cells.Item(26, 11).Value = 2504.5  # K26: 2503 -> 2504.5
$ws.Cells.Item(26, 12).Value = 0  # L26: 2500 -> 0
$ws.Cells.Item(26, 13).Value = -2209.5  # M26: -2208 -> -2209.5
$ws.Cells.Item(26, 14).ClearContents()  # N26: -3090 -> (removed)

$ws.Cells.Item(32, 8).Value = 1130.4  # H32: 933.3333 -> 1130.4
$ws.Cells.Item(32, 9).Value = 1130.4  # I32: 933.3333 -> 1130.4
$ws.Cells.Item(32, 11).Value = 1130.4  # K32: 933.3333 -> 1130.4
$ws.Cells.Item(32, 13).Value = -813.4000000000001  # M32: -616.3333 -> -813.4000000000001

$ws.Cells.Item(42, 8).Value = 0  # H42: 10000 -> 0
$ws.Cells.Item(42, 10).Value = 0  # J42: 10000 -> 0
$ws.Cells.Item(42, 12).Value = 0  # L42: 10000 -> 0
$ws.Cells.Item(42, 14).ClearContents()  # N42: -11126 -> (removed)

$ws.Cells.Item(46, 8).Value = 3630.7693  # H46: 2980.125 -> 3630.7693
$ws.Cells.Item(46, 10).Value = 6285.7144  # J46: 4448.2 -> 6285.7144
$ws.Cells.Item(46, 12).Value = 6285.7144  # L46: 4448.2 -> 6285.7144
$ws.Cells.Item(46, 14).Value = -6661.7144  # N46: -4824.2 -> -6661.7144

$ws.Cells.Item(49, 8).Value = 0  # H49: 10000 -> 0
$ws.Cells.Item(49, 10).Value = 0  # J49: 10000 -> 0
$ws.Cells.Item(49, 12).Value = 0  # L49: 10000 -> 0
$ws.Cells.Item(49, 14).ClearContents()  # N49: -10294 -> (removed)

$ws.Cells.Item(51, 8).Value = 10084  # H51: 8000 -> 10084
$ws.Cells.Item(51, 10).Value = 10084  # J51: 8000 -> 10084
$ws.Cells.Item(51, 12).Value = 10084  # L51: 8000 -> 10084
$ws.Cells.Item(51, 14).Value = -11040  # N51: -8956 -> -11040

$ws.Cells.Item(54, 8).Value = 10084  # H54: 14000 -> 10084
$ws.Cells.Item(54, 10).Value = 10084  # J54: 14000 -> 10084
$ws.Cells.Item(54, 12).Value = 10084  # L54: 14000 -> 10084
$ws.Cells.Item(54, 14).Value = -11372  # N54: -15288 -> -11372

$ws.Cells.Item(55, 8).Value = 975  # H55: 869.1539 -> 975
$ws.Cells.Item(55, 9).Value = 1021.4286  # I55: 770 -> 1021.4286
$ws.Cells.Item(55, 10).Value = 910  # J55: 1199.6666 -> 910
$ws.Cells.Item(55, 11).Value = 1021.4286  # K55: 770 -> 1021.4286
$ws.Cells.Item(55, 12).Value = 910  # L55: 1199.6666 -> 910
$ws.Cells.Item(55, 13).Value = -848.4286  # M55: -597 -> -848.4286
$ws.Cells.Item(55, 14).Value = -1256  # N55: -1545.6666 -> -1256

$ws.Cells.Item(100, 8).Value = 1950.75  # H100: 1601 -> 1950.75
$ws.Cells.Item(100, 9).Value = 1003  # I100: 1301.3334 -> 1003
$ws.Cells.Item(100, 10).Value = 2266.6667  # J100: 2500 -> 2266.6667
$ws.Cells.Item(100, 11).Value = 1003  # K100: 1301.3334 -> 1003
$ws.Cells.Item(100, 12).Value = 2266.6667  # L100: 2500 -> 2266.6667
$ws.Cells.Item(100, 13).Value = -462  # M100: -760.3334 -> -462
$ws.Cells.Item(100, 14).Value = -3348.6667  # N100: -3582 -> -3348.6667

$ws.Cells.Item(122, 8).Value = 22742036  # H122: 25015840 -> 22742036
$ws.Cells.Item(122, 9).Value = 27791066  # I122: 27790956 -> 27791066
$ws.Cells.Item(122, 10).Value = 21402.5  # J122: 39800 -> 21402.5
$ws.Cells.Item(122, 11).Value = 83373198  # K122: 83372868 -> 83373198
$ws.Cells.Item(122, 12).Value = 64207.5  # L122: 119400 -> 64207.5
$ws.Cells.Item(122, 13).Value = -83370748  # M122: -83370418 -> -83370748
$ws.Cells.Item(122, 14).Value = -69107.5  # N122: -124300 -> -69107.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 0  # H24: 200 -> 0
$ws.Cells.Item(24, 9).Value = 0  # I24: 200 -> 0
$ws.Cells.Item(24, 11).Value = 0  # K24: 200 -> 0
$ws.Cells.Item(24, 13).ClearContents()  # M24: 30 -> (removed)

$ws.Cells.Item(29, 8).Value = 0  # H29: 800 -> 0
$ws.Cells.Item(29, 9).Value = 0  # I29: 800 -> 0
$ws.Cells.Item(29, 11).Value = 0  # K29: 800 -> 0
$ws.Cells.Item(29, 13).ClearContents()  # M29: -510 -> (removed)

$ws.Cells.Item(62, 8).Value = 25006066  # H62: 33340966 -> 25006066
$ws.Cells.Item(62, 9).Value = 29416352  # I62: 38467348 -> 29416352
$ws.Cells.Item(62, 10).Value = 14434  # J62: 19501 -> 14434
$ws.Cells.Item(62, 11).Value = 29416352  # K62: 38467348 -> 29416352
$ws.Cells.Item(62, 12).Value = 14434  # L62: 19501 -> 14434
$ws.Cells.Item(62, 13).Value = -29415728  # M62: -38466724 -> -29415728
$ws.Cells.Item(62, 14).Value = -15682  # N62: -20749 -> -15682

$ws.Cells.Item(65, 8).Value = 25006066  # H65: 33340966 -> 25006066
$ws.Cells.Item(65, 9).Value = 29416352  # I65: 38467348 -> 29416352
$ws.Cells.Item(65, 10).Value = 14434  # J65: 19501 -> 14434
$ws.Cells.Item(65, 11).Value = 147081760  # K65: 192336740 -> 147081760
$ws.Cells.Item(65, 12).Value = 72170  # L65: 97505 -> 72170
$ws.Cells.Item(65, 13).Value = -147078640  # M65: -192333620 -> -147078640
$ws.Cells.Item(65, 14).Value = -78410  # N65: -103745 -> -78410
